$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.685.18"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "2.510.97"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "2.508.99"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  +1.40%  "
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("E12").Value = "  +6.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.75%  "
$ws.Range("D14").Value = "2.971.32"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").Value = "69.414.70"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "2.504.62"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("E20").Value = "  -2.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").Value = "0.0₃0892"
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "462.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.51%  "
$ws.Range("E33").Value = "  -5.47%  "
$ws.Range("E34").Value = "  -1.21%  "
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.116"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  -0.39%  "
$ws.Range("E42").Value = "  -1.82%  "
$ws.Range("E43").Value = "  -1.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("E45").Value = "  -5.26%  "
$ws.Range("E46").Value = "  -6.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.76%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.522"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.01%  "
